$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the style of the
# existing header row (e.g. H1) which uses bold font + border + centered top alignment.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in column I (I0) and column J (IF) values for rows 2-8
$iValues = @(1, 1, 1, 1, 1, 2, 6)
$jValues = @(2, 5, 3, 5, 4, 6, 6)

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r - 2]
    $ws.Cells.Item($r, 10).Value = $jValues[$r - 2]
}
